$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DAMSLTag (column I) and DialogAct (column J) values for re-annotated rows
# following a re-run of SGNN dialog act tagging.
$updates = @(
    @{Row=2; DAMSLTag='sv'; DialogAct='Statement-opinion'},
    @{Row=7; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=17; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=22; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=29; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=33; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=36; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=42; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=50; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=55; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=59; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=62; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=67; DAMSLTag='%'; DialogAct='Uninterpretable'},
    @{Row=74; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=76; DAMSLTag='b'; DialogAct='Acknowledge (Backchannel)'},
    @{Row=79; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=89; DAMSLTag='%'; DialogAct='Uninterpretable'},
    @{Row=94; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=98; DAMSLTag='qy'; DialogAct='Yes-No-Question'},
    @{Row=99; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=100; DAMSLTag='ba'; DialogAct='Appreciation'},
    @{Row=103; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=107; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=109; DAMSLTag='ba'; DialogAct='Appreciation'},
    @{Row=115; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=121; DAMSLTag='b'; DialogAct='Acknowledge (Backchannel)'},
    @{Row=123; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=130; DAMSLTag='sv'; DialogAct='Statement-opinion'},
    @{Row=139; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=140; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=145; DAMSLTag='sv'; DialogAct='Statement-opinion'},
    @{Row=149; DAMSLTag='%'; DialogAct='Uninterpretable'},
    @{Row=154; DAMSLTag='%'; DialogAct='Uninterpretable'},
    @{Row=157; DAMSLTag='ba'; DialogAct='Appreciation'},
    @{Row=161; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=162; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=164; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=165; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=166; DAMSLTag='sv'; DialogAct='Statement-opinion'},
    @{Row=168; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=173; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=176; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=184; DAMSLTag='%'; DialogAct='Uninterpretable'},
    @{Row=185; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=186; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=187; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=189; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=190; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=193; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=195; DAMSLTag='sv'; DialogAct='Statement-opinion'},
    @{Row=203; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=218; DAMSLTag='%'; DialogAct='Uninterpretable'},
    @{Row=219; DAMSLTag='qy'; DialogAct='Yes-No-Question'},
    @{Row=224; DAMSLTag='%'; DialogAct='Uninterpretable'},
    @{Row=237; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=257; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=263; DAMSLTag='sv'; DialogAct='Statement-opinion'},
    @{Row=265; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=280; DAMSLTag='%'; DialogAct='Uninterpretable'},
    @{Row=284; DAMSLTag='ba'; DialogAct='Appreciation'},
    @{Row=286; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=293; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=295; DAMSLTag='%'; DialogAct='Uninterpretable'},
    @{Row=306; DAMSLTag='b'; DialogAct='Acknowledge (Backchannel)'},
    @{Row=308; DAMSLTag='sv'; DialogAct='Statement-opinion'},
    @{Row=310; DAMSLTag='sv'; DialogAct='Statement-opinion'},
    @{Row=325; DAMSLTag='b'; DialogAct='Acknowledge (Backchannel)'},
    @{Row=327; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=329; DAMSLTag='ba'; DialogAct='Appreciation'},
    @{Row=335; DAMSLTag='b'; DialogAct='Acknowledge (Backchannel)'},
    @{Row=345; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=355; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=368; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=375; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=380; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=386; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=389; DAMSLTag='sv'; DialogAct='Statement-opinion'},
    @{Row=406; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=409; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=424; DAMSLTag='%'; DialogAct='Uninterpretable'},
    @{Row=458; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=468; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=470; DAMSLTag='%'; DialogAct='Uninterpretable'},
    @{Row=489; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=504; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=528; DAMSLTag='%'; DialogAct='Uninterpretable'},
    @{Row=532; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=538; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=546; DAMSLTag='b'; DialogAct='Acknowledge (Backchannel)'},
    @{Row=553; DAMSLTag='sv'; DialogAct='Statement-opinion'},
    @{Row=572; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=578; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}

Write-Output "Updated $($updates.Count) rows"